$d = $word.ActiveDocument

# Merge the three runs "<id>" + "p072r_1" + "</id>" into a single run
# with the text "<id>p072r_1</id>", keeping the formatting of the
# first ("<id>") run. A literal Find/Replace over the whole span
# achieves exactly this merge in Word.
$d.Content.Find.Execute("<id>p072r_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p072r_1</id>", 2)

# Same merge for the second occurrence, "<id>" + "p072r_2" + "</id>".
$d.Content.Find.Execute("<id>p072r_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p072r_2</id>", 2)
